$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("metro_budget")

# --- Section around row 64-70: XLOOKUP lookups for a handful of departments ---
$ws.Range("B65").Formula = '=XLOOKUP(A65,A:A,D:D)'
$ws.Range("C65").Formula = '=XLOOKUP(A65,A:A,I:I)'
$ws.Range("D65").Formula = '=XLOOKUP(A65,A:A,N:N)'

$ws.Range("B66:B70").Formula = '=XLOOKUP(A66,A:A,D:D)'
$ws.Range("C66:C70").Formula = '=XLOOKUP(A66,A:A,I:I)'
$ws.Range("D66:D70").Formula = '=XLOOKUP(A66,A:A,N:N)'

# --- Section around row 73-79: INDEX/MATCH lookups for the same departments ---
$ws.Range("B74").Formula = '=INDEX(D:D,MATCH(A74,A:A,0))'
$ws.Range("C74").Formula = '=INDEX(I:I,MATCH(A74,A:A,0))'
$ws.Range("D74").Formula = '=INDEX(I:I,MATCH(A74,A:A,0))'

$ws.Range("B75:B79").Formula = '=INDEX(D:D,MATCH(A75,A:A,0))'
$ws.Range("C75:C79").Formula = '=INDEX(I:I,MATCH(A75,A:A,0))'
$ws.Range("D75:D79").Formula = '=INDEX(I:I,MATCH(A75,A:A,0))'

# --- Section around row 81-86: pick a department in B82 and look up its values ---
$ws.Range("B82").Value = "Codes Administration"

$ws.Range("B84").Formula = '=INDEX($B$2:$B$52,MATCH($B$82,$A$2:$A$52,0))'
$ws.Range("C84").Formula = '=INDEX($C$2:$C$52,MATCH($B$82,$A$2:$A$52,0))'
$ws.Range("B85").Formula = '=INDEX($G$2:$G$52,MATCH($B$82,$A$2:$A$52,0))'
$ws.Range("C85").Formula = '=INDEX($H$2:$H$52,MATCH($B$82,$A$2:$A$52,0))'
$ws.Range("B86").Formula = '=INDEX($L$2:$L$52,MATCH($B$82,$A$2:$A$52,0))'
$ws.Range("C86").Formula = '=INDEX($M$2:$M$52,MATCH($B$82,$A$2:$A$52,0))'

# --- Data validation: B82 gets its own list validation, split off of the old B82:B83 blank rule ---
$ws.Range("B82").Validation.Delete()
$ws.Range("B82").Validation.Add(3, 1, 1, '$A2:$A$52')

# --- Selection / view state ---
$ws.Range("H92").Select() | Out-Null
